{"js": "// Word JS API (Office.js) edit script.\n// Applies the three text corrections described by the commit\n// \"corre\u00e7\u00e3o do segundo feedback\":\n//\n// 1) \"...dos pedidos feitos. Facilitando a troca...\"\n//      -> \"...dos pedidos feitos, facilitando a troca...\"\n// 2) \"O cliente poder\u00e1 controlar os gastos ... estabelecimento, e os ganhos...\"\n//      -> \"O gerente poder\u00e1 controlar os gastos ... estabelecimento e os ganhos...\"\n// 3) \"De acordo com o tempo as cores do background onde o gerente tem cada\n//      tarefa especificada, v\u00e3o mudar de mais fortes (vermelho) a mais\n//      fracas (azul) de acordo com o tempo de limite cr\u00edtico ao est\u00e1vel\"\n//      -> \"De acordo com as cores do background onde o gerente tem cada\n//      tarefa especificada, cores mais fortes (vermelho) identificam\n//      tarefas pendentes e cores mais fracas (azul) tarefas finalizadas de\n//      acordo com o tempo de limite\"\n\nasync function replaceOnce(body, searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + searchText);\n  }\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// Change 1\nawait replaceOnce(\n  body,\n  \"A numera\u00e7\u00e3o das mesas garantir\u00e1 um maior controle dos pedidos feitos. Facilitando a troca de informa\u00e7\u00f5es entre as \u00e1reas da lanchonete e a cobran\u00e7a pelo consumo.\",\n  \"A numera\u00e7\u00e3o das mesas garantir\u00e1 um maior controle dos pedidos feitos, facilitando a troca de informa\u00e7\u00f5es entre as \u00e1reas da lanchonete e a cobran\u00e7a pelo consumo.\"\n);\n\n// Change 2\nawait replaceOnce(\n  body,\n  \"O cliente poder\u00e1 controlar os gastos referentes \u00e0s despesas para manuten\u00e7\u00e3o do estabelecimento, e os ganhos relativos aos recebidos pela presta\u00e7\u00e3o de servi\u00e7o.\",\n  \"O gerente poder\u00e1 controlar os gastos referentes \u00e0s despesas para manuten\u00e7\u00e3o do estabelecimento e os ganhos relativos aos recebidos pela presta\u00e7\u00e3o de servi\u00e7o.\"\n);\n\n// Change 3\nawait replaceOnce(\n  body,\n  \"De acordo com o tempo as cores do background onde o gerente tem cada tarefa especificada, v\u00e3o mudar de mais fortes (vermelho) a mais fracas (azul) de acordo com o tempo de limite cr\u00edtico ao est\u00e1vel\",\n  \"De acordo com as cores do background onde o gerente tem cada tarefa especificada, cores mais fortes (vermelho) identificam tarefas pendentes e cores mais fracas (azul) tarefas finalizadas de acordo com o tempo de limite\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the three text corrections described by the commit\n# \"corre\u00e7\u00e3o do segundo feedback\":\n#\n# 1) \"...dos pedidos feitos. Facilitando a troca...\"\n#      -> \"...dos pedidos feitos, facilitando a troca...\"\n# 2) \"O cliente poder\u00e1 controlar os gastos ... estabelecimento, e os ganhos...\"\n#      -> \"O gerente poder\u00e1 controlar os gastos ... estabelecimento e os ganhos...\"\n# 3) \"De acordo com o tempo as cores do background onde o gerente tem cada\n#      tarefa especificada, v\u00e3o mudar de mais fortes (vermelho) a mais\n#      fracas (azul) de acordo com o tempo de limite cr\u00edtico ao est\u00e1vel\"\n#      -> \"De acordo com as cores do background onde o gerente tem cada\n#      tarefa especificada, cores mais fortes (vermelho) identificam\n#      tarefas pendentes e cores mais fracas (azul) tarefas finalizadas de\n#      acordo com o tempo de limite\"\n\n$d = $word.ActiveDocument\n\n# wdReplace constants\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $found = $find.Execute($findText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# Change 1\nReplace-Text \"A numera\u00e7\u00e3o das mesas garantir\u00e1 um maior controle dos pedidos feitos. Facilitando a troca de informa\u00e7\u00f5es entre as \u00e1reas da lanchonete e a cobran\u00e7a pelo consumo.\" \"A numera\u00e7\u00e3o das mesas garantir\u00e1 um maior controle dos pedidos feitos, facilitando a troca de informa\u00e7\u00f5es entre as \u00e1reas da lanchonete e a cobran\u00e7a pelo consumo.\"\n\n# Change 2\nReplace-Text \"O cliente poder\u00e1 controlar os gastos referentes \u00e0s despesas para manuten\u00e7\u00e3o do estabelecimento, e os ganhos relativos aos recebidos pela presta\u00e7\u00e3o de servi\u00e7o.\" \"O gerente poder\u00e1 controlar os gastos referentes \u00e0s despesas para manuten\u00e7\u00e3o do estabelecimento e os ganhos relativos aos recebidos pela presta\u00e7\u00e3o de servi\u00e7o.\"\n\n# Change 3\nReplace-Text \"De acordo com o tempo as cores do background onde o gerente tem cada tarefa especificada, v\u00e3o mudar de mais fortes (vermelho) a mais fracas (azul) de acordo com o tempo de limite cr\u00edtico ao est\u00e1vel\" \"De acordo com as cores do background onde o gerente tem cada tarefa especificada, cores mais fortes (vermelho) identificam tarefas pendentes e cores mais fracas (azul) tarefas finalizadas de acordo com o tempo de limite\"\n"}
